$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.923.07"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.630.59"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.82"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.34"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.861.51"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.625.76"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.57"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "27.918.70"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.10"
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.26"
$ws.Range("E23").Value = "  -5.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.72"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.08"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").Value = "1.398.88"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  +12.26%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.556"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.863"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.40"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.84"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.48"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "1.771.71"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.08"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("E51").Value = "  -0.27%  "
